$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.785703420639038
$ws.Range("B1").Value = 4.090403079986572
$ws.Range("C1").Value = 2.074463844299316
$ws.Range("D1").Value = 1.603772878646851
$ws.Range("E1").Value = 1.450293779373169
